$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B3").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B4").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B5").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B6").Value = "ΔCFI = 0; ΔRMSEA = 0.01."
$ws.Range("B7").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B8").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B9").Value = "ΔCFI = 0; ΔRMSEA = 0."
$ws.Range("B10").Value = "ΔCFI = 0; ΔRMSEA = 0."
